# "Update countries & provincias Spain"
# - Senegal was missing/mis-ordered relative to Mauricio: the country row that
#   used to read "Mauricio" now shows "Senegal" (with Senegal's refreshed
#   numbers) and the row below now shows "Mauricio" (keeping the figures that
#   previously belonged to the old "Mauricio" row).
# - A handful of per-country case/death counters were refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Senegal / Mauricio rows (row 111 <-> row 112 labels) ---
$ws.Range("A111").Value = "Senegal"
$ws.Range("A112").Value = "Mauricio"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 644823
$ws.Range("C4").Value = 734
$ws.Range("E4").Value = 567533
$ws.Range("G4").Value = 51
$ws.Range("H4").Value = 28580

# --- Row 8: Alemania ---
$ws.Range("B8").Value = 135549
$ws.Range("C8").Value = 796
$ws.Range("E8").Value = 54699

# --- Row 20: Austria ---
$ws.Range("B20").Value = 14451
$ws.Range("C20").Value = 101
$ws.Range("D20").Value = 8986
$ws.Range("E20").Value = 5072

# --- Row 60: Moldavia ---
$ws.Range("E60").Value = 1761
$ws.Range("G60").Value = 7
$ws.Range("H60").Value = 53

# --- Row 111: now Senegal (updated figures) ---
$ws.Range("B111").Value = 335
$ws.Range("C111").Value = 21
$ws.Range("D111").Value = 194
$ws.Range("E111").Value = 139
$ws.Range("F111").Value = 1
$ws.Range("H111").Value = 2

# --- Row 112: now Mauricio (figures that used to be on the Mauricio row) ---
$ws.Range("B112").Value = 324
$ws.Range("D112").Value = 65
$ws.Range("E112").Value = 250
$ws.Range("F112").Value = 3
$ws.Range("H112").Value = 9

# --- Row 143: Togo ---
$ws.Range("D143").Value = 45
$ws.Range("E143").Value = 31
$ws.Range("G143").Value = 2
$ws.Range("H143").Value = 5
